# Reorder the ranking table on "Foglio1": move the 13-row VALUE block
# (currently rows 47-59) up to the top of the data (rows 2-14), pushing
# the GROWTH / MOMENTUM / QUALITY blocks (currently rows 2-46) down by
# 13 rows so they occupy rows 15-59. Header row 1 is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Read the whole data block (58 rows x 4 cols) as a single COM array.
$full = $ws.Range("A2:D59").Value()

$rowCount = 58
$colCount = 4
$blockSize = 13          # number of VALUE rows being moved to the front
$restSize = $rowCount - $blockSize

$newArr = New-Object 'object[,]' $rowCount, $colCount

# VALUE block (old rows 47-59 => array rows 46-58, 1-based) goes first.
for ($i = 0; $i -lt $blockSize; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $newArr[$i, $j] = $full[($restSize + 1 + $i), (1 + $j)]
    }
}

# Remaining rows (old rows 2-46 => array rows 1-45, 1-based) follow.
for ($i = 0; $i -lt $restSize; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $newArr[($blockSize + $i), $j] = $full[(1 + $i), (1 + $j)]
    }
}

$ws.Range("A2:D59").Value = $newArr

# Update the visible selection to match the reordered VALUE block.
$ws.Range("A2:D14").Select()
